$wb = $excel.ActiveWorkbook
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# --- Overview sheet: row 3 is the "b.md" row. Mark it ready for handoff. ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-25 09:00:24"
$wsOverview.Range("D3").NumberFormat = $dateFmt

# --- zh-cn sheet: row 3 is the "b.md" row. New handoff file generated. ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-25 09:00:19"
$wsZhCn.Range("E3").NumberFormat = $dateFmt
foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 is the "b.md" row. New handoff file generated. ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-25 09:00:24"
$wsDeDe.Range("E3").NumberFormat = $dateFmt
foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
